$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.329.46'
$ws.Range("E2").Value = '  +0.56%  '
$ws.Range("D3").Value = '1.874.08'
$ws.Range("E3").Value = '  +0.67%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7118'
$ws.Range("E5").Value = '  +0.44%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '241.95'
$ws.Range("E6").Value = '  +0.40%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07796'
$ws.Range("E8").Value = '  +1.96%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3110'
$ws.Range("E9").Value = '  +0.42%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '25.17'
$ws.Range("E10").Value = '  +1.94%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08403'
$ws.Range("E11").Value = '  +0.57%  '
$ws.Range("D12").Value = '1.874.03'
$ws.Range("E12").Value = '  +0.59%  '
$ws.Range("E13").Value = '  +1.11%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7111'
$ws.Range("E14").Value = '  +0.36%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.18'
$ws.Range("E15").Value = '  +0.05%  '
$ws.Range("D16").Value = '29.339.67'
$ws.Range("E16").Value = '  +0.59%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.087'
$ws.Range("E17").Value = '  +2.88%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008203'
$ws.Range("E18").Value = '  +5.16%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '240.15'
$ws.Range("E19").Value = '  -1.13%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.20'
$ws.Range("E20").Value = '  +0.89%  '
$ws.Range("D21").Value = '2.121.91'
$ws.Range("E21").Value = '  +0.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.762'
$ws.Range("E23").Value = '  -1.38%  '
$ws.Range("E24").Value = '  +0.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1599'
$ws.Range("E25").Value = '  +1.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.67'
$ws.Range("E26").Value = '  -0.37%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.029'
$ws.Range("E27").Value = '  +0.89%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.51'
$ws.Range("E28").Value = '  +0.42%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.509'
$ws.Range("E29").Value = '  +0.75%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.423'
$ws.Range("E30").Value = '  +0.46%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.291'
$ws.Range("E31").Value = '  -2.94%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.315'
$ws.Range("E32").Value = '  +1.90%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05288'
$ws.Range("E33").Value = '  +2.84%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.940'
$ws.Range("E35").Value = '  +1.19%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7480'
$ws.Range("E36").Value = '  -5.82%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.701'
$ws.Range("E37").Value = '  +0.69%  '
$ws.Range("E38").Value = '  +1.47%  '
$ws.Range("D39").Value = '1.227.13'
$ws.Range("E39").Value = '  +5.17%  '
$ws.Range("E40").Value = '  +1.22%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.538'
$ws.Range("E41").Value = '  +5.44%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '110.40'
$ws.Range("E42").Value = '  +8.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8871'
$ws.Range("E43").Value = '  -0.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '72.65'
$ws.Range("E44").Value = '  -0.37%  '
$ws.Range("E45").Value = '  +0.11%  '
$ws.Range("D46").Value = '2.019.97'
$ws.Range("E46").Value = '  +0.57%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.802'
$ws.Range("E47").Value = '  +1.63%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5193'
$ws.Range("E48").Value = '  -0.11%  '
$ws.Range("E49").Value = '  +2.56%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.394'
$ws.Range("E50").Value = '  +0.65%  '
$ws.Range("E51").Value = '  +1.03%  '
